# Regenerate orders with updated distance/size codes.
#
# The experiment's "Distance" codes change D64->D69, D51->D55, D80->D86,
# and the "Size" code S30 becomes S31 (S20 / S25 are unchanged). These
# tokens appear embedded inside many text values across the sheet
# (Condition, Filename_Left, Filename_Right, Distance, Size columns), so
# walk every used cell and rewrite any string value containing one of the
# old tokens.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value()
        if ($val -is [string]) {
            $newVal = $val.Replace("D64", "D69").Replace("D51", "D55").Replace("D80", "D86").Replace("S30", "S31")
            if ($newVal -ne $val) {
                $cell.Value = $newVal
            }
        }
    }
}
